$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I55").Value = "sv"
$ws.Range("J55").Value = "Statement-opinion"
$ws.Range("I59").Value = "sd"
$ws.Range("J59").Value = "Statement-non-opinion"
$ws.Range("I75").Value = "b"
$ws.Range("J75").Value = "Acknowledge (Backchannel)"
$ws.Range("I84").Value = "aa"
$ws.Range("J84").Value = "Agree/Accept"
$ws.Range("I91").Value = "aa"
$ws.Range("J91").Value = "Agree/Accept"
$ws.Range("I103").Value = "sd"
$ws.Range("J103").Value = "Statement-non-opinion"
$ws.Range("I112").Value = "sd"
$ws.Range("J112").Value = "Statement-non-opinion"
$ws.Range("I115").Value = "aa"
$ws.Range("J115").Value = "Agree/Accept"
$ws.Range("I116").Value = "aa"
$ws.Range("J116").Value = "Agree/Accept"
$ws.Range("I127").Value = "sd"
$ws.Range("J127").Value = "Statement-non-opinion"
$ws.Range("I137").Value = "sd"
$ws.Range("J137").Value = "Statement-non-opinion"
$ws.Range("I139").Value = "sd"
$ws.Range("J139").Value = "Statement-non-opinion"
$ws.Range("I141").Value = "sd"
$ws.Range("J141").Value = "Statement-non-opinion"
$ws.Range("I148").Value = "b"
$ws.Range("J148").Value = "Acknowledge (Backchannel)"
$ws.Range("I162").Value = "aa"
$ws.Range("J162").Value = "Agree/Accept"
$ws.Range("I185").Value = "sd"
$ws.Range("J185").Value = "Statement-non-opinion"
$ws.Range("I194").Value = "sv"
$ws.Range("J194").Value = "Statement-opinion"
$ws.Range("I196").Value = "sv"
$ws.Range("J196").Value = "Statement-opinion"
$ws.Range("I201").Value = "sd"
$ws.Range("J201").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "ba"
$ws.Range("J214").Value = "Appreciation"
$ws.Range("I215").Value = "sd"
$ws.Range("J215").Value = "Statement-non-opinion"
$ws.Range("I235").Value = "qy"
$ws.Range("J235").Value = "Yes-No-Question"
$ws.Range("I240").Value = "ba"
$ws.Range("J240").Value = "Appreciation"
$ws.Range("I248").Value = "sd"
$ws.Range("J248").Value = "Statement-non-opinion"
$ws.Range("I253").Value = "sd"
$ws.Range("J253").Value = "Statement-non-opinion"
$ws.Range("I263").Value = "aa"
$ws.Range("J263").Value = "Agree/Accept"
$ws.Range("I299").Value = "b"
$ws.Range("J299").Value = "Acknowledge (Backchannel)"
$ws.Range("I315").Value = "aa"
$ws.Range("J315").Value = "Agree/Accept"
$ws.Range("I327").Value = "sv"
$ws.Range("J327").Value = "Statement-opinion"
$ws.Range("I329").Value = "sd"
$ws.Range("J329").Value = "Statement-non-opinion"
$ws.Range("I349").Value = "sd"
$ws.Range("J349").Value = "Statement-non-opinion"
$ws.Range("I360").Value = "sd"
$ws.Range("J360").Value = "Statement-non-opinion"
$ws.Range("I366").Value = "sd"
$ws.Range("J366").Value = "Statement-non-opinion"
$ws.Range("I367").Value = "aa"
$ws.Range("J367").Value = "Agree/Accept"
$ws.Range("I368").Value = "sd"
$ws.Range("J368").Value = "Statement-non-opinion"
